$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "JsonField": add a "WarnMsg" row after "CheckMsg", drop the old
# "FitRate" row and append three new rows (DrawdownDate / MaturityDate /
# FirstAdjRateDate) at the bottom of the table.
# ---------------------------------------------------------------------------
$wsJson = $wb.Worksheets.Item("JsonField")

# Insert a new row 6 (pushes RateIncr / IndividualIncr / FitRate down by one)
# copying the formatting of row 5 ("CheckMsg") so borders/fonts match.
$wsJson.Range("A6:H6").Insert(-4121)
$wsJson.Range("A5:H5").Copy()
$wsJson.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsJson.Cells.Item(6, 2).Value = 2
$wsJson.Cells.Item(6, 3).Value = "WarnMsg"
$wsJson.Cells.Item(6, 4).Value = "提示訊息"
$wsJson.Cells.Item(6, 5).Value = "nVARCHAR2"
$wsJson.Cells.Item(6, 6).Value = $null
$wsJson.Cells.Item(6, 7).Value = $null
$wsJson.Cells.Item(6, 8).Value = $null

# Renumber the two rows that shifted down (content unchanged, only SEQ++)
$wsJson.Cells.Item(7, 2).Value = 3
$wsJson.Cells.Item(8, 2).Value = 4

# The former "FitRate" row is now row 9 - turn it into the "DrawdownDate" row
$wsJson.Cells.Item(9, 2).Value = 5
$wsJson.Cells.Item(9, 3).Value = "DrawdownDate"
$wsJson.Cells.Item(9, 4).Value = "撥款日期"
$wsJson.Cells.Item(9, 5).Value = "DECIMALD"
$wsJson.Cells.Item(9, 6).Value = 8
$wsJson.Cells.Item(9, 7).Value = $null
$wsJson.Cells.Item(9, 8).Value = $null

# Append two more rows (MaturityDate / FirstAdjRateDate), copying row 9's
# formatting.
$wsJson.Range("A9:H9").Copy()
$wsJson.Range("A10:H11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsJson.Cells.Item(10, 2).Value = 6
$wsJson.Cells.Item(10, 3).Value = "MaturityDate"
$wsJson.Cells.Item(10, 4).Value = "到期日"
$wsJson.Cells.Item(10, 5).Value = "DECIMALD"
$wsJson.Cells.Item(10, 6).Value = 8
$wsJson.Cells.Item(10, 7).Value = $null
$wsJson.Cells.Item(10, 8).Value = $null

$wsJson.Cells.Item(11, 2).Value = 7
$wsJson.Cells.Item(11, 3).Value = "FirstAdjRateDate"
$wsJson.Cells.Item(11, 4).Value = "首次調整日期"
$wsJson.Cells.Item(11, 5).Value = "DECIMALD"
$wsJson.Cells.Item(11, 6).Value = 8
$wsJson.Cells.Item(11, 7).Value = $null
$wsJson.Cells.Item(11, 8).Value = $null

# ---------------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved state of the workbook.
# ---------------------------------------------------------------------------
$wsJson.Range("D17").Select()

$wsDbd = $wb.Worksheets.Item("DBD")
$wsDbd.Activate()
$wsDbd.Application.ActiveWindow.ScrollRow = 28
$wsDbd.Range("G34").Select()
